# Apply the "cryptos list" price/volume refresh described in the commit.
# Values that look like plain numbers (e.g. "25.80", "162.00") must stay as
# literal text -- same as the source data -- so we briefly force a text
# number format before writing them, then clear formatting again so the
# cell keeps the workbook default style (matches the original, unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.775.09'
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").Value = '2.497.09'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '587.53'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.65%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '176.28'
$c.ClearFormats()
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.78%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.142'
$c.ClearFormats()
$ws.Range("E9").Value = '  +3.85%  '
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").Value = '2.958.24'
$ws.Range("E13").Value = '  +1.09%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.80'
$c.ClearFormats()
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("D15").Value = '67.752.08'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("D17").Value = '2.503.63'
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("E18").Value = '  +0.97%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.53'
$c.ClearFormats()
$ws.Range("E19").Value = '  +1.34%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '351.16'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.21%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.10'
$c.ClearFormats()
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  +0.04%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '70.86'
$c.ClearFormats()
$ws.Range("E24").Value = '  +2.58%  '
$ws.Range("E25").Value = '  -0.18%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '9.17'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  +1.04%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("D29").Value = '0.0₃0908'
$ws.Range("E29").Value = '  +0.97%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '507.86'
$c.ClearFormats()
$ws.Range("E30").Value = '  +0.00%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.83'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("E32").Value = '  +3.26%  '
$ws.Range("E33").Value = '  +1.01%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  +5.83%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '162.00'
$c.ClearFormats()
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("E37").Value = '  +0.06%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '18.41'
$c.ClearFormats()
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  +3.61%  '
$ws.Range("E42").Value = '  +1.15%  '
$ws.Range("E43").Value = '  +1.15%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.ClearFormats()
$ws.Range("E44").Value = '  +2.46%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '145.44'
$c.ClearFormats()
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("E46").Value = '  +2.42%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.516'
$c.ClearFormats()
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0744'
$c.ClearFormats()
$ws.Range("E48").Value = '  +2.19%  '
$ws.Range("B49").Value = 'Optimism'
$ws.Range("C49").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.59'
$c.ClearFormats()
$ws.Range("E49").Value = '  +1.72%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.587'
$c.ClearFormats()
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.17'
$c.ClearFormats()
$ws.Range("E51").Value = '  +0.26%  '
